# SCD0010-115 - Penyesuaian Kriteria Account Plan SRM
# Rename sheet, update the TC_ID text "DGS-182" -> "SCD0010-115",
# resize column B to fit the new text, and move the active selection
# to reflect where the author last clicked (O3, scrolled to show col H).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab (SCD0167 -> SCD0010)
$ws.Name = "SCD0010"

# Update the TC_ID values in rows 2 and 3 from "DGS-182" to "SCD0010-115"
$ws.Range("B2").Value = "SCD0010-115"
$ws.Range("B3").Value = "SCD0010-115"

# Widen column B so the longer TC_ID text fits (was 9 chars wide).
# Target stored width is ~12.42578125; this engine quantizes ColumnWidth
# onto a coarse internal grid, so 11.6 is the closest input that lands
# on the nearest achievable stored width (12.5).
$ws.Columns.Item(2).ColumnWidth = 11.6

# Move the visible selection/active cell to O3 (scrolled so column H is leftmost)
$null = $ws.Range("O3").Select()
